$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 658635213881.7474
$ws.Range("C3").Value = 213648463155.1176
$ws.Range("C4").Value = 84291615051.44429
$ws.Range("C5").Value = 34149224073.46932
$ws.Range("C6").Value = 29763354750.16835
$ws.Range("C7").Value = 25147785370.52747
$ws.Range("C8").Value = 12810329657.90201
$ws.Range("C9").Value = 9795719791.685951
$ws.Range("C10").Value = 9346057965.164358
$ws.Range("C11").Value = 8190902529.1807
$ws.Range("C12").Value = 7335984498.418167

# Row 13 becomes Polygon / MATIC-USD
$ws.Range("A13").Value = "Polygon"
$ws.Range("B13").Value = "MATIC-USD"
$ws.Range("C13").Value = 5842193209.048368

# Row 14 becomes Chainlink / LINK-USD
$ws.Range("A14").Value = "Chainlink"
$ws.Range("B14").Value = "LINK-USD"
$ws.Range("C14").Value = 5555046571.631498

# Row 15 stays Wrapped Bitcoin / WBTC-USD, value changes
$ws.Range("C15").Value = 5502141881.324093

# Row 16 becomes Dai / DAI-USD
$ws.Range("A16").Value = "Dai"
$ws.Range("B16").Value = "DAI-USD"
$ws.Range("C16").Value = 5346266803.592254
